$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "22.458.44"
Set-TextValue $ws.Range("E2") "  +0.55%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.571.39"
Set-TextValue $ws.Range("E3") "  +0.39%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  -0.39%  "

# Row 5 - USDC
Set-TextValue $ws.Range("E5") "  -0.29%  "

# Row 6 - BNB
Set-TextValue $ws.Range("D6") "290.36"
Set-TextValue $ws.Range("E6") "  +0.27%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.3692"
Set-TextValue $ws.Range("E7") "  -1.79%  "

# Row 8 - OKB
Set-TextValue $ws.Range("D8") "49.94"
Set-TextValue $ws.Range("E8") "  +1.52%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.3382"
Set-TextValue $ws.Range("E9") "  +0.60%  "

# Row 10 - Polygon
Set-TextValue $ws.Range("D10") "1.150"
Set-TextValue $ws.Range("E10") "  +2.27%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.07555"
Set-TextValue $ws.Range("E11") "  +0.64%  "

# Row 12 - BinanceUSD
Set-TextValue $ws.Range("D12") "1.001"
Set-TextValue $ws.Range("E12") "  -0.41%  "

# Row 13 - Solana
Set-TextValue $ws.Range("D13") "21.22"
Set-TextValue $ws.Range("E13") "  +2.10%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("E14") "  +2.31%  "

# Row 15 - Chainlink
Set-TextValue $ws.Range("D15") "6.991"
Set-TextValue $ws.Range("E15") "  +1.96%  "

# Row 16 - WrappedEther
Set-TextValue $ws.Range("D16") "1.570.13"
Set-TextValue $ws.Range("E16") "  +0.43%  "

# Row 17 - ShibaInu
Set-TextValue $ws.Range("D17") "0.00001123"
Set-TextValue $ws.Range("E17") "  +0.80%  "

# Row 18 - Litecoin
Set-TextValue $ws.Range("D18") "90.43"
Set-TextValue $ws.Range("E18") "  +1.32%  "

# Row 19 - TRON
Set-TextValue $ws.Range("D19") "0.06774"
Set-TextValue $ws.Range("E19") "  +0.97%  "

# Row 20 - Dai
Set-TextValue $ws.Range("D20") "1.001"
Set-TextValue $ws.Range("E20") "  -0.38%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "6.361"
Set-TextValue $ws.Range("E21") "  +3.22%  "

# Row 22 - Avalanche
Set-TextValue $ws.Range("D22") "16.42"
Set-TextValue $ws.Range("E22") "  +0.58%  "

# Row 23 - Cosmos
Set-TextValue $ws.Range("D23") "12.21"
Set-TextValue $ws.Range("E23") "  +3.17%  "

# Row 24 - WrappedBTC
Set-TextValue $ws.Range("D24") "22.464.16"
Set-TextValue $ws.Range("E24") "  +0.61%  "

# Row 25 - Toncoin
Set-TextValue $ws.Range("D25") "2.354"
Set-TextValue $ws.Range("E25") "  -1.01%  "

# Row 26 - LidoDAOToken
Set-TextValue $ws.Range("D26") "2.674"
Set-TextValue $ws.Range("E26") "  +0.56%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "20.03"
Set-TextValue $ws.Range("E27") "  +0.19%  "

# Row 28 - Monero
Set-TextValue $ws.Range("D28") "149.23"
Set-TextValue $ws.Range("E28") "  +1.21%  "

# Row 29 - HuobiToken
Set-TextValue $ws.Range("E29") "  +1.42%  "

# Row 30 - BitcoinCash
Set-TextValue $ws.Range("D30") "124.96"
Set-TextValue $ws.Range("E30") "  -0.12%  "

# Row 31 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D31") "1.746.54"
Set-TextValue $ws.Range("E31") "  +0.42%  "

# Row 32 - ImmutableX
Set-TextValue $ws.Range("E32") "  +8.69%  "

# Row 33 - Filecoin
Set-TextValue $ws.Range("D33") "6.241"
Set-TextValue $ws.Range("E33") "  +5.20%  "

# Row 34 - WEMIXTOKEN
Set-TextValue $ws.Range("D34") "2.017"
Set-TextValue $ws.Range("E34") "  +0.03%  "

# Row 35 - FraxShare
Set-TextValue $ws.Range("D35") "9.812"
Set-TextValue $ws.Range("E35") "  -0.17%  "

# Row 36 - Stellar
Set-TextValue $ws.Range("D36") "0.08380"
Set-TextValue $ws.Range("E36") "  -0.72%  "

# Row 37 - VeChain
Set-TextValue $ws.Range("D37") "0.02480"
Set-TextValue $ws.Range("E37") "  +1.33%  "

# Row 38 - TrustWalletToken
Set-TextValue $ws.Range("D38") "1.353"
Set-TextValue $ws.Range("E38") "  -3.78%  "

# Row 39 - Algorand
Set-TextValue $ws.Range("D39") "0.2302"
Set-TextValue $ws.Range("E39") "  +1.88%  "

# Row 40 - Hedera
Set-TextValue $ws.Range("D40") "0.06568"
Set-TextValue $ws.Range("E40") "  +2.85%  "

# Row 41 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D41") "5.429"
Set-TextValue $ws.Range("E41") "  +1.73%  "

# Row 42 - Aptos
Set-TextValue $ws.Range("D42") "11.32"
Set-TextValue $ws.Range("E42") "  +3.32%  "

# Row 43 - TheSandbox
Set-TextValue $ws.Range("D43") "0.6257"
Set-TextValue $ws.Range("E43") "  +0.47%  "

# Row 44 - Frax -> EnergySwap (row identity swap)
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D44") "14.14"
Set-TextValue $ws.Range("E44") "  +2.10%  "

# Row 45 - EnergySwap -> Frax (row identity swap)
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D45") "1.001"
Set-TextValue $ws.Range("E45") "  -0.38%  "

# Row 46 - PancakeSwap
Set-TextValue $ws.Range("D46") "3.800"

# Row 47 - Decentraland
Set-TextValue $ws.Range("D47") "0.5878"
Set-TextValue $ws.Range("E47") "  +1.31%  "

# Row 48 - NEARProtocol
Set-TextValue $ws.Range("D48") "2.075"
Set-TextValue $ws.Range("E48") "  +1.52%  "

# Row 49 - Quant
Set-TextValue $ws.Range("D49") "127.97"
Set-TextValue $ws.Range("E49") "  +3.10%  "

# Row 50 - EOS
Set-TextValue $ws.Range("D50") "1.245"
Set-TextValue $ws.Range("E50") "  -0.16%  "

# Row 51 - Cronos
Set-TextValue $ws.Range("D51") "0.07307"
Set-TextValue $ws.Range("E51") "  -0.02%  "
